$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 'Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '26.573.49', '  +0.16%  '),
    @(3, 'Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '1.813.62', '  -0.03%  '),
    @(4, 'TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '1.003', '  -0.43%  '),
    @(5, 'USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '1.003', '  -0.44%  '),
    @(6, 'BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '306.01', '  -0.87%  '),
    @(7, 'XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.4551', '  -0.38%  '),
    @(8, 'Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.3596', '  -2.00%  '),
    @(9, 'OKB', 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb', '46.39', '  +2.48%  '),
    @(10, 'Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.07123', '  -0.12%  '),
    @(11, 'Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '0.8935', '  +1.53%  '),
    @(12, 'TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.07716', '  -0.49%  '),
    @(13, 'Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '19.35', '  -0.06%  '),
    @(14, 'WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '1.804.63', '  -1.73%  '),
    @(15, 'Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '5.258', '  -0.73%  '),
    @(16, 'Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '6.298', '  -1.21%  '),
    @(17, 'Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '86.62', '  -0.12%  '),
    @(18, 'BinanceUSD', 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd', '1.004', '  -0.48%  '),
    @(19, 'ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.000008554', '  -0.52%  '),
    @(20, 'Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '1.002', '  -0.43%  '),
    @(21, 'WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '26.593.56', '  -0.02%  '),
    @(22, 'Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '14.16', '  -0.67%  '),
    @(23, 'Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '4.960', '  -1.11%  '),
    @(24, 'Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '10.55', '  +0.51%  '),
    @(25, 'Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '1.922', '  -3.01%  '),
    @(26, 'Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '151.77', '  +0.22%  '),
    @(27, 'EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '17.82', '  -0.70%  '),
    @(28, 'LidoDAOToken', 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo', '2.023', '  -2.45%  '),
    @(29, 'BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '112.40', '  -0.51%  '),
    @(30, 'InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '4.834', '  -0.43%  '),
    @(31, 'Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.08722', '  +0.22%  '),
    @(32, 'HuobiToken', 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht', '3.117', '  +2.71%  '),
    @(33, 'ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '0.7421', '  +1.31%  '),
    @(34, 'Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '4.427', '  -1.70%  '),
    @(35, 'RenderToken', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr', '2.710', '  +1.12%  '),
    @(36, 'ARBITRUM', 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb', '1.111', '  -0.85%  '),
    @(37, 'TrustWalletToken', 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt', '1.070', '  -1.42%  '),
    @(38, 'VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.01938', '  -0.93%  '),
    @(39, 'MXToken', 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx', '2.918', '  +1.00%  '),
    @(40, 'Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.05079', '  -0.91%  '),
    @(41, 'TheSandbox', 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand', '0.5099', '  +1.96%  '),
    @(42, 'FraxShare', 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs', '6.785', '  -3.04%  '),
    @(43, 'Algorand', 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo', '0.1511', '  -3.07%  '),
    @(44, 'Aptos', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt', '8.028', '  -1.65%  '),
    @(45, 'Decentraland', 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana', '0.4692', '  +1.88%  '),
    @(46, 'PaxDollar', 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp', '1.003', '  -0.51%  '),
    @(47, 'EnergySwap', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens', '10.03', '  +0.89%  '),
    @(48, 'Quant', 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt', '99.26', '  -2.13%  '),
    @(49, 'NEARProtocol', 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near', '1.568', '  -1.45%  '),
    @(50, 'Cronos', 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro', '0.05995', '  -0.09%  '),
    @(51, 'Aave', 'https://coinranking.com/coin/ixgUfzmLR+aave-aave', '63.83', '  -1.03%  ')
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = "'" + $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}
